$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell D5: "spezielle Zutaten wissen" -> "speziellen Zutaten kennen" ---
$ws.Cells.Item(5, 4).Value = "Als ältere Generation möchte man die speziellen Zutaten kennen, um die Rezepte weitergeben zu können"

# --- Fill in previously blank row 10 ---
$ws.Cells.Item(10, 1).Value = "Benutzer"
$ws.Cells.Item(10, 2).Value = "Anrecht"
$ws.Cells.Item(10, 3).Value = "auf genaue Angaben vom Rezept"
$ws.Cells.Item(10, 4).Value = "Als Benutzer muss man genaue Angaben verfügbar haben, um das Rezept exakt nachkochen zu können"
$ws.Cells.Item(10, 4).WrapText = $true

# --- Fill in previously blank row 11 ---
$ws.Cells.Item(11, 1).Value = "ältere Generation"
$ws.Cells.Item(11, 2).Value = "Interesse"
$ws.Cells.Item(11, 3).Value = "einfach bedienbares System"
$ws.Cells.Item(11, 4).Value = "Als ältere Generation möchte man ein einfach bedienbares System haben, um die Rezepte `nproblemlos eintragen zu können"
$ws.Cells.Item(11, 4).WrapText = $true

# --- Row heights: rows 1-13 become 70 (were previously 60 for 1-10, 40 for 11-13) ---
$ws.Rows("1:13").RowHeight = 70

# --- Update sheet view: scroll position and active selection ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("D11").Select()
